# Correccion de observaciones #1
# Adds the two new columns "TIP_UBICACION" (J) and "TIP_CAJA" (K) to the
# Hoja1 table header row, matches the resulting used range / selection,
# and sizes the two new columns the way the author left them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (plain, unstyled text cells - same as the rest of the
# "extra" trailing columns G/H/I which also carry no explicit style).
$ws.Range("J1").Value = "TIP_UBICACION"
$ws.Range("K1").Value = "TIP_CAJA"

# Author resized the two brand-new columns by hand.
$ws.Columns.Item(10).ColumnWidth = 16
$ws.Columns.Item(11).ColumnWidth = 14.86

# Leave the selection where the author left it when they saved the file.
[void]$ws.Range("K3").Select()
